# Applies the cryptos-list refresh described in the commit message.
# Updates Price (D) and Volume(1h) (E) columns for most rows, and for three
# row-pairs (12/13, 21/22, 46/47) the ranking order flips: Coin (B) and
# Link (C) are swapped along with their new Price/Volume figures.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '29.917.54'
$ws.Range("E2").Value = '  +0.15%  '

# Row 3
$ws.Range("D3").Value = '1.875.46'
$ws.Range("E3").Value = '  -0.64%  '

# Row 4
$ws.Range("D4").Value = '''1.002'
$ws.Range("E4").Value = '  +0.12%  '

# Row 5
$ws.Range("D5").Value = '''0.7388'
$ws.Range("E5").Value = '  -3.85%  '

# Row 6
$ws.Range("D6").Value = '''242.66'
$ws.Range("E6").Value = '  +0.03%  '

# Row 7
$ws.Range("E7").Value = '  +0.07%  '

# Row 8
$ws.Range("D8").Value = '''0.3158'
$ws.Range("E8").Value = '  +1.13%  '

# Row 9
$ws.Range("D9").Value = '''0.07201'
$ws.Range("E9").Value = '  +0.52%  '

# Row 10
$ws.Range("E10").Value = '  -4.10%  '

# Row 11
$ws.Range("D11").Value = '''0.08344'
$ws.Range("E11").Value = '  -2.65%  '

# Row 12
$ws.Range("B12").Value = 'WrappedEther'
$ws.Range("C12").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D12").Value = '1.938.62'
$ws.Range("E12").Value = '  +3.33%  '

# Row 13
$ws.Range("B13").Value = 'Polygon'
$ws.Range("C13").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D13").Value = '''0.7495'
$ws.Range("E13").Value = '  -1.79%  '

# Row 14
$ws.Range("D14").Value = '''5.389'
$ws.Range("E14").Value = '  +0.47%  '

# Row 15
$ws.Range("D15").Value = '''92.31'
$ws.Range("E15").Value = '  -1.36%  '

# Row 16
$ws.Range("D16").Value = '29.924.03'
$ws.Range("E16").Value = '  +0.30%  '

# Row 17
$ws.Range("D17").Value = '''6.090'
$ws.Range("E17").Value = '  -0.98%  '

# Row 18
$ws.Range("D18").Value = '''247.72'
$ws.Range("E18").Value = '  +1.35%  '

# Row 19
$ws.Range("E19").Value = '  -1.59%  '

# Row 20
$ws.Range("D20").Value = '''0.000007829'
$ws.Range("E20").Value = '  +0.38%  '

# Row 21
$ws.Range("B21").Value = 'Dai'
$ws.Range("C21").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D21").Value = '''1.001'
$ws.Range("E21").Value = '  +0.16%  '

# Row 22
$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").Value = '2.141.56'
$ws.Range("E22").Value = '  +1.01%  '

# Row 23
$ws.Range("D23").Value = '''8.018'
$ws.Range("E23").Value = '  +0.23%  '

# Row 24
$ws.Range("D24").Value = '''1.002'
$ws.Range("E24").Value = '  +0.11%  '

# Row 25
$ws.Range("D25").Value = '''0.1552'
$ws.Range("E25").Value = '  -4.98%  '

# Row 26
$ws.Range("D26").Value = '''9.273'
$ws.Range("E26").Value = '  -1.17%  '

# Row 27
$ws.Range("D27").Value = '''164.56'
$ws.Range("E27").Value = '  +1.08%  '

# Row 28
$ws.Range("E28").Value = '  -0.37%  '

# Row 29
$ws.Range("E29").Value = '  -0.07%  '

# Row 30
$ws.Range("D30").Value = '''1.508'
$ws.Range("E30").Value = '  +3.18%  '

# Row 31
$ws.Range("E31").Value = '  +1.86%  '

# Row 32
$ws.Range("D32").Value = '''1.535'
$ws.Range("E32").Value = '  -0.26%  '

# Row 33
$ws.Range("D33").Value = '''4.268'
$ws.Range("E33").Value = '  +4.30%  '

# Row 34
$ws.Range("D34").Value = '''0.05315'
$ws.Range("E34").Value = '  -2.53%  '

# Row 35
$ws.Range("D35").Value = '''1.233'
$ws.Range("E35").Value = '  -0.51%  '

# Row 36
$ws.Range("D36").Value = '''0.7484'
$ws.Range("E36").Value = '  +0.87%  '

# Row 37
$ws.Range("D37").Value = '''0.9992'
$ws.Range("E37").Value = '  -0.16%  '

# Row 38
$ws.Range("D38").Value = '''2.697'
$ws.Range("E38").Value = '  -0.11%  '

# Row 39
$ws.Range("D39").Value = '''0.01965'
$ws.Range("E39").Value = '  +0.54%  '

# Row 40
$ws.Range("E40").Value = '  -0.96%  '

# Row 41
$ws.Range("D41").Value = '''0.4543'
$ws.Range("E41").Value = '  +1.73%  '

# Row 42
$ws.Range("D42").Value = '1.104.55'
$ws.Range("E42").Value = '  -0.38%  '

# Row 43
$ws.Range("D43").Value = '''6.122'
$ws.Range("E43").Value = '  +0.87%  '

# Row 44
$ws.Range("D44").Value = '''72.35'
$ws.Range("E44").Value = '  -0.87%  '

# Row 45
$ws.Range("D45").Value = '''0.8613'
$ws.Range("E45").Value = '  +1.10%  '

# Row 46
$ws.Range("B46").Value = 'Quant'
$ws.Range("C46").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D46").Value = '''104.26'
$ws.Range("E46").Value = '  +1.84%  '

# Row 47
$ws.Range("B47").Value = 'PaxDollar'
$ws.Range("C47").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D47").Value = '''1.002'
$ws.Range("E47").Value = '  +0.17%  '

# Row 48
$ws.Range("D48").Value = '''1.855'
$ws.Range("E48").Value = '  -0.18%  '

# Row 49
$ws.Range("D49").Value = '''7.604'
$ws.Range("E49").Value = '  -0.59%  '

# Row 50
$ws.Range("D50").Value = '''9.515'
$ws.Range("E50").Value = '  -2.70%  '

# Row 51
$ws.Range("D51").Value = '2.037.98'
$ws.Range("E51").Value = '  -0.85%  '
